$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.064850662639623
$ws.Range("D2").Value = 1.067496491402062
$ws.Range("E2").Value = 1.069373432351804
$ws.Range("F2").Value = 1.078631132388632
$ws.Range("I2").Value = 1.040942000987361
$ws.Range("J2").Value = 1.069808788005184
$ws.Range("K2").Value = 1.070204197677716
$ws.Range("L2").Value = 1.072076112540631
$ws.Range("M2").Value = 1.081309298727051
$ws.Range("N2").Value = 1.071328039192935
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.066690141083253
$ws.Range("D3").Value = 1.06925374974298
$ws.Range("E3").Value = 1.071058848150994
$ws.Range("F3").Value = 1.080508516253083
$ws.Range("I3").Value = 1.041318684544749
$ws.Range("J3").Value = 1.071299427179064
$ws.Range("K3").Value = 1.071774910914318
$ws.Range("L3").Value = 1.073575531931038
$ws.Range("M3").Value = 1.083002023489604
$ws.Range("N3").Value = 1.07282079524542
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.067876923370908
$ws.Range("D4").Value = 1.070387654884151
$ws.Range("E4").Value = 1.072146411900874
$ws.Range("F4").Value = 1.081720454628794
$ws.Range("I4").Value = 1.041559635488666
$ws.Range("J4").Value = 1.072260206516456
$ws.Range("K4").Value = 1.07278765814295
$ws.Range("L4").Value = 1.074542276373316
$ws.Range("M4").Value = 1.084094042076517
$ws.Range("N4").Value = 1.073782938999674
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.068375030861913
$ws.Range("D5").Value = 1.070863608954309
$ws.Range("E5").Value = 1.072602918080968
$ws.Range("F5").Value = 1.082229287562106
$ws.Range("I5").Value = 1.041660268209831
$ws.Range("J5").Value = 1.072663231883842
$ws.Range("K5").Value = 1.073212568360176
$ws.Range("L5").Value = 1.07494787722822
$ws.Range("M5").Value = 1.08455235584127
$ws.Range("N5").Value = 1.07418653670931
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.068458617986612
$ws.Range("D6").Value = 1.070943480839869
$ws.Range("E6").Value = 1.072679526552387
$ws.Range("F6").Value = 1.08231468430127
$ws.Range("I6").Value = 1.041677126134903
$ws.Range("J6").Value = 1.07273085007296
$ws.Range("K6").Value = 1.073283863329125
$ws.Range("L6").Value = 1.075015931821617
$ws.Range("M6").Value = 1.084629264003469
$ws.Range("N6").Value = 1.074254250924011
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.067883582288373
$ws.Range("D7").Value = 1.070394017489546
$ws.Range("E7").Value = 1.07215251450999
$ws.Range("F7").Value = 1.081727256272372
$ws.Range("I7").Value = 1.041560982747159
$ws.Range("J7").Value = 1.072265595223345
$ws.Range("K7").Value = 1.07279333912801
$ws.Range("L7").Value = 1.07454769922798
$ws.Range("M7").Value = 1.084100169094303
$ws.Range("N7").Value = 1.073788335359145
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.065473055183416
$ws.Range("D8").Value = 1.068091029057123
$ws.Range("E8").Value = 1.069943659858505
$ws.Range("F8").Value = 1.079266204637007
$ws.Range("I8").Value = 1.041069883072758
$ws.Range("J8").Value = 1.070313344791904
$ws.Range("K8").Value = 1.070735784175842
$ws.Range("L8").Value = 1.072583577796671
$ws.Range("M8").Value = 1.081882053481569
$ws.Range("N8").Value = 1.071833312508165
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.061197903181685
$ws.Range("D9").Value = 1.064007933972077
$ws.Range("E9").Value = 1.066027591769108
$ws.Range("F9").Value = 1.074906869038383
$ws.Range("I9").Value = 1.040182941664778
$ws.Range("J9").Value = 1.066843734370953
$ws.Range("K9").Value = 1.067081773952819
$ws.Range("L9").Value = 1.069095229634702
$ws.Range("M9").Value = 1.077947541161576
$ws.Range("N9").Value = 1.068358774842446
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.058328166271141
$ws.Range("D10").Value = 1.061268058135115
$ws.Range("E10").Value = 1.063399892739871
$ws.Range("F10").Value = 1.071984313050985
$ws.Range("I10").Value = 1.039576862496118
$ws.Range("J10").Value = 1.064509888528414
$ws.Range("K10").Value = 1.06462574136022
$ws.Range("L10").Value = 1.066750368730747
$ws.Range("M10").Value = 1.075306085550962
$ws.Range("N10").Value = 1.066021614671115
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.057080623735565
$ws.Range("D11").Value = 1.060077202008666
$ws.Range("E11").Value = 1.062257817722204
$ws.Range("F11").Value = 1.070714692710395
$ws.Range("I11").Value = 1.039310851767806
$ws.Range("J11").Value = 1.063494169602831
$ws.Range("K11").Value = 1.063557283767971
$ws.Range("L11").Value = 1.065730232878404
$ws.Range("M11").Value = 1.074157700347805
$ws.Range("N11").Value = 1.065004453308153
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.056616468275974
$ws.Range("D12").Value = 1.059634172826555
$ws.Range("E12").Value = 1.061832940234463
$ws.Range("F12").Value = 1.070242456603882
$ws.Range("I12").Value = 1.039211500836138
$ws.Range("J12").Value = 1.063116094584592
$ws.Range("K12").Value = 1.063159644006111
$ws.Range("L12").Value = 1.065350570419075
$ws.Range("M12").Value = 1.073730425258274
$ws.Range("N12").Value = 1.064625841380022
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.056716066068846
$ws.Range("D13").Value = 1.059729235742833
$ws.Range("E13").Value = 1.061924108076568
$ws.Range("F13").Value = 1.070343782277231
$ws.Range("I13").Value = 1.039232836598759
$ws.Range("J13").Value = 1.063197229085879
$ws.Range("K13").Value = 1.063244974101982
$ws.Range("L13").Value = 1.065432043020814
$ws.Range("M13").Value = 1.073822109834811
$ws.Range("N13").Value = 1.064707091101607
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.057042272187696
$ws.Range("D14").Value = 1.060040595317511
$ws.Range("E14").Value = 1.062222710792701
$ws.Range("F14").Value = 1.070675670760059
$ws.Range("I14").Value = 1.039302650492937
$ws.Range("J14").Value = 1.063462934062849
$ws.Range("K14").Value = 1.06352443049914
$ws.Range("L14").Value = 1.065698865033462
$ws.Range("M14").Value = 1.074122396331172
$ws.Range("N14").Value = 1.064973173410122
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.057243156817385
$ws.Range("D15").Value = 1.060232341891756
$ws.Range("E15").Value = 1.062406601744709
$ws.Range("F15").Value = 1.070880072413284
$ws.Range("I15").Value = 1.039345593052508
$ws.Range("J15").Value = 1.063626538197413
$ws.Range("K15").Value = 1.063696510597101
$ws.Range("L15").Value = 1.065863164466969
$ws.Range("M15").Value = 1.074307317579246
$ws.Range("N15").Value = 1.065137009881322
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.058410855845552
$ws.Range("D16").Value = 1.061346995386921
$ws.Range("E16").Value = 1.063475597013306
$ws.Range("F16").Value = 1.072068484641882
$ws.Range("I16").Value = 1.039594440959308
$ws.Range("J16").Value = 1.064577188422283
$ws.Range("K16").Value = 1.064696544851485
$ws.Range("L16").Value = 1.0668179692103
$ws.Range("M16").Value = 1.075382201111053
$ws.Range("N16").Value = 1.066089010138553
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.059141987698176
$ws.Range("D17").Value = 1.062044976557629
$ws.Range("E17").Value = 1.064144994321554
$ws.Range("F17").Value = 1.072812822306517
$ws.Range("I17").Value = 1.039749575708371
$ws.Range("J17").Value = 1.065172114998484
$ws.Range("K17").Value = 1.065322493518261
$ws.Range("L17").Value = 1.067415596377774
$ws.Range("M17").Value = 1.076055196902772
$ws.Range("N17").Value = 1.066684781578732
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.059567970625458
$ws.Range("D18").Value = 1.062451667088986
$ws.Range("E18").Value = 1.064535032352253
$ws.Range("F18").Value = 1.073246584434365
$ws.Range("I18").Value = 1.039839718688348
$ws.Range("J18").Value = 1.065518630120225
$ws.Range("K18").Value = 1.065687119650825
$ws.Range("L18").Value = 1.067763720721075
$ws.Range("M18").Value = 1.076447299802346
$ws.Range("N18").Value = 1.06703178879169
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.059713140082322
$ws.Range("D19").Value = 1.062590265713854
$ws.Range("E19").Value = 1.064667956264165
$ws.Range("F19").Value = 1.073394419342926
$ws.Range("I19").Value = 1.039870396853463
$ws.Range("J19").Value = 1.065636699329716
$ws.Range("K19").Value = 1.065811367039678
$ws.Range("L19").Value = 1.067882344396095
$ws.Range("M19").Value = 1.076580922021395
$ws.Range("N19").Value = 1.067150025673002
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.059063593349603
$ws.Range("D20").Value = 1.061970134415241
$ws.Range("E20").Value = 1.064073216883003
$ws.Range("F20").Value = 1.072733003183043
$ws.Range("I20").Value = 1.039732966891629
$ws.Range("J20").Value = 1.065108336389543
$ws.Range("K20").Value = 1.065255384846064
$ws.Range("L20").Value = 1.067351524501324
$ws.Range("M20").Value = 1.075983036908957
$ws.Range("N20").Value = 1.066620912396851
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.056946233848679
$ws.Range("D21").Value = 1.059948926955917
$ws.Range("E21").Value = 1.062134798110526
$ws.Range("F21").Value = 1.070577955785996
$ws.Range("I21").Value = 1.039282107068456
$ws.Range("J21").Value = 1.063384712575043
$ws.Range("K21").Value = 1.063442158882963
$ws.Range("L21").Value = 1.065620313136081
$ws.Range("M21").Value = 1.074033989305494
$ws.Range("N21").Value = 1.064894840838831
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.055610542377159
$ws.Range("D22").Value = 1.058674098510497
$ws.Range("E22").Value = 1.060912209108938
$ws.Range("F22").Value = 1.069219263347762
$ws.Range("I22").Value = 1.0389954908183
$ws.Range("J22").Value = 1.062296411606735
$ws.Range("K22").Value = 1.062297664331325
$ws.Range("L22").Value = 1.064527549495797
$ws.Range("M22").Value = 1.07280440643338
$ws.Range("N22").Value = 1.063804994358319
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.056319044519826
$ws.Range("D23").Value = 1.059350296638557
$ws.Range("E23").Value = 1.061560696066012
$ws.Range("F23").Value = 1.06993989257742
$ws.Range("I23").Value = 1.039147731332795
$ws.Range("J23").Value = 1.062873781821901
$ws.Range("K23").Value = 1.062904810505747
$ws.Range("L23").Value = 1.065107256196512
$ws.Range("M23").Value = 1.073456630667845
$ws.Range("N23").Value = 1.064383184505415
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.059099017845147
$ws.Range("D24").Value = 1.06200395368635
$ws.Range("E24").Value = 1.064105651284651
$ws.Range("F24").Value = 1.072769071236797
$ws.Range("I24").Value = 1.039740472765369
$ws.Range("J24").Value = 1.065137156728248
$ws.Range("K24").Value = 1.065285709846482
$ws.Range("L24").Value = 1.067380477249906
$ws.Range("M24").Value = 1.076015644276092
$ws.Range("N24").Value = 1.066649773663742
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.062306508251329
$ws.Range("D25").Value = 1.065066573924347
$ws.Range("E25").Value = 1.067042909553867
$ws.Range("F25").Value = 1.076036658645669
$ws.Range("I25").Value = 1.040414821828592
$ws.Range("J25").Value = 1.067744303385439
$ws.Range("K25").Value = 1.068029880683072
$ws.Range("L25").Value = 1.070000383314932
$ws.Range("M25").Value = 1.078967882298434
$ws.Range("N25").Value = 1.06926062276823
